$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.474.37"
$ws.Range("E2").Value = "  -0.50%  "
$ws.Range("D3").Value = "1.823.63"
$ws.Range("E3").Value = "  -0.94%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.56"
$ws.Range("E5").Value = "  -0.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4237"
$ws.Range("E7").Value = "  -0.53%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3629"
$ws.Range("E8").Value = "  +0.48%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07191"
$ws.Range("E9").Value = "  -1.51%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8586"
$ws.Range("E10").Value = "  -1.92%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.58"
$ws.Range("E11").Value = "  -0.05%  "
$ws.Range("D12").Value = "1.766.43"
$ws.Range("E12").Value = "  -3.26%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.406"
$ws.Range("E13").Value = "  +1.53%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.475"
$ws.Range("E14").Value = "  -0.31%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06935"
$ws.Range("E15").Value = "  -0.25%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.005"
$ws.Range("E16").Value = "  +0.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "80.15"
$ws.Range("E17").Value = "  +1.18%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008861"
$ws.Range("E18").Value = "  -0.77%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.002"
$ws.Range("E19").Value = "  +0.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.33"
$ws.Range("E20").Value = "  +0.10%  "
$ws.Range("D21").Value = "27.299.40"
$ws.Range("E21").Value = "  -0.95%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.119"
$ws.Range("E22").Value = "  +3.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.94"
$ws.Range("E23").Value = "  +6.18%  "
$ws.Range("D24").Value = "2.028.08"
$ws.Range("E24").Value = "  -2.72%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.982"
$ws.Range("E25").Value = "  -0.25%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.80"
$ws.Range("E26").Value = "  -0.46%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.71"
$ws.Range("E27").Value = "  +1.23%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.121"
$ws.Range("E28").Value = "  -1.37%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "113.99"
$ws.Range("E29").Value = "  -4.21%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.801"
$ws.Range("E30").Value = "  -3.88%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08837"
$ws.Range("E31").Value = "  -0.47%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.977"
$ws.Range("E32").Value = "  +1.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7433"
$ws.Range("E33").Value = "  -1.96%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.530"
$ws.Range("E34").Value = "  +0.84%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.121"
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.002"
$ws.Range("E36").Value = "  +0.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.087"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05277"
$ws.Range("E38").Value = "  -2.73%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01926"
$ws.Range("E39").Value = "  +0.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.785"
$ws.Range("E40").Value = "  -1.33%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5043"
$ws.Range("E41").Value = "  -0.18%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1641"
$ws.Range("E42").Value = "  -1.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.467"
$ws.Range("E43").Value = "  -0.90%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.296"
$ws.Range("E44").Value = "  -0.82%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.35"
$ws.Range("E45").Value = "  -0.10%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "105.38"
$ws.Range("E46").Value = "  -0.42%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.06450"
$ws.Range("E47").Value = "  -1.48%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4664"
$ws.Range("E48").Value = "  +0.95%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.001"
$ws.Range("E49").Value = "  +0.08%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.613"
$ws.Range("E50").Value = "  -1.38%  "

# Row 51 special changes (coin renamed from Aave to RenderToken)
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.720"
$ws.Range("E51").Value = "  -0.44%  "
